# Update spawning temp estimate approach:
#  - column F header renamed "prop.ripe" -> "percent.ripe"
#  - column F values recomputed as a percentage (ripe / (mature+ripe+spent) * 100)
#    instead of a raw proportion, using full double precision
#  - column F formatted with a 2-decimal numeric format
#  - column F width set to fit the new values/header
#  - selection moved to F8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text for column F
$ws.Range("F1").Value = "percent.ripe"

# Recompute every data row's percent-ripe value from the source counts
# (mature = C, ripe = D, spent = E) at full floating point precision.
$lastRow = 277
for ($r = 2; $r -le $lastRow; $r++) {
    $mature = $ws.Cells.Item($r, 3).Value2
    $ripe   = $ws.Cells.Item($r, 4).Value2
    $spent  = $ws.Cells.Item($r, 5).Value2

    $denom = $mature + $ripe + $spent
    if ($denom -eq 0) {
        $percentRipe = 0
    } else {
        $percentRipe = $ripe / $denom * 100
    }

    $ws.Cells.Item($r, 6).Value = $percentRipe
}

# Apply a 2-decimal-place number format to the recomputed column (creates the
# new cellXfs entry and tags every F2:F277 cell with it).
$ws.Range("F2:F277").NumberFormat = "0.00"

# Resize column F to fit its new contents/header.
$ws.Columns.Item(6).ColumnWidth = 12.666666666667

# Move the active selection to F8.
$ws.Range("F8").Select()
